$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns F, G, H with same style as existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style from an existing header cell (E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean outlier flag values for rows 2-12 in columns F (KNN), G (SVM), H (RF)
$knnVals = @(0,0,1,0,0,0,1,0,0,1,0)
$svmVals = @(0,0,0,0,0,0,0,0,0,0,0)
$rfVals  = @(0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt 11; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = [bool]($knnVals[$i])
    $ws.Cells.Item($row, 7).Value = [bool]($svmVals[$i])
    $ws.Cells.Item($row, 8).Value = [bool]($rfVals[$i])
}
